# "use case 3.1 pronto"
# Adds a new "UC 3.1" (Listar Insumos) sheet, cloned from "UC 1.1" (Listar
# Clientes), tweaks a few wording fixes on both sheets, and restores the
# various cell-selection / active-sheet bookmarks that Excel re-wrote the
# last time the workbook was saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Wording fixes on "UC 1.1" (Listar Clientes)
# ---------------------------------------------------------------------
$uc11 = $wb.Worksheets.Item("UC 1.1")
$uc11.Range("C5").Value2 = "PR2 - O Funcionário ou Administrador deve estar na tela de busca de clientes."
$uc11.Range("C7").Value2 = "AF1 - Se não houver nenhum cliete com alguma ligação com o que foi pesquisado o sistema retornará uma mensagemde erro;"
$uc11.Range("C9").Value2 = "Aparecerá uma lista com os clientes com o que foram descritos na pesquisa"

# ---------------------------------------------------------------------
# 2. Create "UC 3.1" (Listar Insumos) as a copy of "UC 1.1", placed right
#    after the last sheet ("UC 3.0").
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$uc11.Copy([System.Reflection.Missing]::Value, $lastSheet)
$uc31 = $wb.Worksheets.Item($wb.Worksheets.Count)
$uc31.Name = "UC 3.1"

# Re-word the copied cells for the "Insumos" (supplies) flavour of the use case
$uc31.Range("C2").Value2 = "Listar Insumos"
$uc31.Range("C5").Value2 = "PR2 - O Funcionário ou Administrador deve estar na tela de busca de insumos."
$uc31.Range("C7").Value2 = "AF1 - Se não houver nenhum Insumo com alguma ligação com o que foi pesquisado o sistema retornará uma mensagemde erro;"
$uc31.Range("C9").Value2 = "Aparecerá uma lista com os Insumos  que foram descritos na pesquisa"

# Match the original column widths used on the new sheet
$uc31.Columns.Item(2).ColumnWidth = 17.7
$uc31.Columns.Item(3).ColumnWidth = 112.7

# Selection left on the new sheet
$uc31.Range("C9").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Restore the lingering cell selections on the other sheets
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Product Backlog")
$backlog.Range("B3").Select() | Out-Null

$uc10 = $wb.Worksheets.Item("UC 1.0")
$uc10.Range("C14").Select() | Out-Null

$uc11.Range("B2:C10").Select() | Out-Null
$uc11.Range("B2").Activate() | Out-Null

$uc30 = $wb.Worksheets.Item("UC 3.0")
$uc30.Range("B2:C14").Select() | Out-Null
$uc30.Range("C14").Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3

# ---------------------------------------------------------------------
# 4. Finally, leave "UC 3.1" as the active sheet/tab (matches the diff's
#    activeTab bump + tabSelected moving to the new sheet).
# ---------------------------------------------------------------------
$uc31.Activate() | Out-Null
$uc31.Range("C9").Select() | Out-Null
